# Sample Project / Main.xlsx - "Rules" sheet
#
# Change cell B11 (the "R40" rule-id label in the decision table) to the
# literal text "1". The cell must keep its original formatting (style index
# / General number format) and the new value must be stored as text, not as
# a number - exactly as when a user types a numeric-looking value into a
# General cell in Excel using a leading apostrophe ('1) to force text entry.
#
# A direct `Range.Value = "1"` would have Excel auto-detect the numeric
# string and store it as a Number (and would also silently reformat the
# cell if we pre-set NumberFormat="@" to force text). To avoid touching the
# cell's existing style we stage the text value in a scratch cell, then use
# PasteSpecial (values only) to land just the value onto B11, leaving its
# formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$target = $ws.Range("B11")

# Stage the text value "1" (apostrophe forces text, not number) in an
# unused scratch cell far away from the used range.
$scratch = $ws.Cells.Item(50, 50)
$scratch.Value = "'1"

# Copy only the value onto B11 so its existing formatting/style is kept.
$scratch.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false

# Clean up the scratch cell.
$scratch.Clear()

$wb.Save()
